$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain-looking decimal numbers as text
# (t="inlineStr" in the sheet XML). Assigning a numeric-looking string
# straight to .Value would let Excel auto-convert it to a real number,
# which changes its representation (e.g. trailing zeros get dropped, or
# floating point noise like 94.510000000000005 appears). Forcing the
# Text number format on those specific cells first keeps the new value
# stored as text, matching the original cell formatting.
$textCells = @("D5", "D6", "D10", "D11", "D16", "D17", "D22", "D23", "D29", "D31", "D32", "D33", "D34", "D39", "D42", "D46", "D47", "D48", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '44.610.56'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.246.69'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').Value = '306.14'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '94.51'
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').Value = '34.64'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').Value = '0.0802'
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.590.65'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.237.94'
$ws.Range('E15').Value = '  -3.87%  '
$ws.Range('D16').Value = '0.832'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').Value = '13.56'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '44.371.24'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').Value = '0.0₃0938'
$ws.Range('E19').Value = '  -2.56%  '
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('E21').Value = '  -2.72%  '
$ws.Range('D22').Value = '65.22'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '237.22'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  +5.17%  '
$ws.Range('E28').Value = '  -1.53%  '
$ws.Range('D29').Value = '36.93'
$ws.Range('E29').Value = '  -2.86%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').Value = '19.93'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').Value = '148.47'
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = '0.0782'
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').Value = '2.62'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('E36').Value = '  +1.94%  '
$ws.Range('E37').Value = '  -1.21%  '
$ws.Range('E38').Value = '  +5.64%  '
$ws.Range('D39').Value = '15.20'
$ws.Range('E39').Value = '  +6.51%  '
$ws.Range('E40').Value = '  -5.62%  '
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('D42').Value = '0.0298'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '1.804.17'
$ws.Range('E44').Value = '  +3.03%  '
$ws.Range('E45').Value = '  +10.76%  '
$ws.Range('D46').Value = '81.73'
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('D47').Value = '0.187'
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('D48').Value = '98.37'
$ws.Range('E48').Value = '  -1.39%  '
$ws.Range('E49').Value = '  -2.04%  '
$ws.Range('D50').Value = '68.32'
$ws.Range('E50').Value = '  +2.56%  '
$ws.Range('D51').Value = '53.75'
$ws.Range('E51').Value = '  -1.33%  '
